# Weekly update: insert a new most-recent week's price row for Jengibre
# (Mercado Mayorista Lo Valledor de Santiago) at row 10, pushing the
# existing historical rows (old rows 10-63) down by one to rows 11-64.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 10; this shifts rows 10..63 -> 11..64
# and carries the row-9/row-10 cell formatting (incl. the date style on
# column D) down onto the newly inserted row, matching the original
# workbook's per-row style pattern.
$ws.Rows(10).Insert()

# Populate the newly inserted row 10 with this week's observation.
$ws.Cells.Item(10, 1).Value = 6
$ws.Cells.Item(10, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(10, 3).Value = "Metropolitana"
$ws.Cells.Item(10, 4).Value = 44602
$ws.Cells.Item(10, 5).Value = 13
$ws.Cells.Item(10, 6).Value = 100114007
$ws.Cells.Item(10, 7).Value = "Jengibre"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 200
$ws.Cells.Item(10, 11).Value = 12000
$ws.Cells.Item(10, 12).Value = 13000
$ws.Cells.Item(10, 13).Value = 12400
$ws.Cells.Item(10, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(10, 15).Value = "Perú"
$ws.Cells.Item(10, 16).Value = 954
$ws.Cells.Item(10, 17).Value = 13
$ws.Cells.Item(10, 18).Value = "Hortaliza"
